$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 9: A9 = "3ColumnModel", B9 = "B"
$ws.Range("A9").Value = "3ColumnModel"
$ws.Range("B9").Value = "B"

# Update the selected cell on the sheet to reflect new selection (B11)
$ws.Range("B11").Select()
